$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footprint values to explicitly note "imperial" sizing
$ws.Range("D2").Value = "0603 (imperial)"
$ws.Range("D3").Value = "0603 (imperial)"
$ws.Range("D4").Value = "1206 (imperial)"
$ws.Range("D5").Value = "0603 (imperial)"
$ws.Range("D8").Value = "0603 (imperial)"

# Update the active cell selection
$ws.Range("D8").Select()
